# hw2/doubleor.xlsx -- "prob 1 mostly working"
#
# - 5qubit sheet: tweak the C2 coefficient, mark the rows that hit the
#   (tied) minimum Objective value with an "x" in column I.
# - 6qubit sheet: fill in the A2:U2 coefficients (H column formulas then
#   recompute automatically).
# - Leave 6qubit as the active/selected tab, with new active-cell
#   selections on both sheets.

$wb = $excel.ActiveWorkbook

$ws5 = $wb.Worksheets.Item("5qubit")
$ws6 = $wb.Worksheets.Item("6qubit")

# ---------------------------------------------------------------------
# 5qubit ("5qubit" == sheet1)
# ---------------------------------------------------------------------

# Coefficient tweak: C2 -1 -> -2 (H5:H36 formulas recompute automatically)
$ws5.Range("C2").Value = -2

# Mark every row whose recomputed Objective (column H) tied for the
# minimum (-2) with an "x" in column I.
$ws5.Range("I9").Value = "x"
$ws5.Range("I11").Value = "x"
$ws5.Range("I14").Value = "x"
$ws5.Range("I15").Value = "x"
$ws5.Range("I22").Value = "x"
$ws5.Range("I23").Value = "x"
$ws5.Range("I30").Value = "x"
$ws5.Range("I31").Value = "x"

# ---------------------------------------------------------------------
# 6qubit (sheet2)
# ---------------------------------------------------------------------

# Fill in the A2:U2 model coefficients; H5:H68 formulas recompute
# automatically off of these.
$ws6.Range("A2").Value = -1
$ws6.Range("B2").Value = -1
$ws6.Range("C2").Value = 0.5
$ws6.Range("D2").Value = -0.5
$ws6.Range("E2").Value = -1
$ws6.Range("F2").Value = -1
$ws6.Range("G2").Value = 1
$ws6.Range("H2").Value = 2
$ws6.Range("I2").Value = 0
$ws6.Range("J2").Value = 0
$ws6.Range("K2").Value = 0
$ws6.Range("L2").Value = 0
$ws6.Range("M2").Value = 2
$ws6.Range("N2").Value = 0
$ws6.Range("O2").Value = 0
$ws6.Range("P2").Value = -2
$ws6.Range("Q2").Value = 0
$ws6.Range("R2").Value = 0
$ws6.Range("S2").Value = 1
$ws6.Range("T2").Value = 2
$ws6.Range("U2").Value = 2

# ---------------------------------------------------------------------
# View state: 6qubit becomes the active/selected sheet, each sheet gets
# a new active cell.
# ---------------------------------------------------------------------

$ws5.Range("P2").Select()
$ws6.Range("V2").Select()

$ws6.Select()
$ws6.Range("V2").Select()
